$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 72 (the "Xenium_v1_skin_human" / human-skin-preview-data dataset row)
$ws.Rows.Item(72).Delete()
